# Camel case naming for examples and exercises.
#
# Each target slide's title placeholder (Shapes.Item(1)) holds two runs,
# e.g. run1="Example: Token " + run2="Tester". The edit joins the two
# words into a single camelCase word that lives entirely in the second
# run, e.g. run1="Example: " + run2="TokenTester". We do this with a
# single Characters(start, length) assignment that spans from the start
# of the first word through the end of the title text - PowerPoint's
# TextRange keeps that replaced span inside the (pre-existing) second
# run's paragraph position, which collapses the run split back down to
# exactly two runs: the untouched "Example: " prefix, and the new
# camelCase word.

$p = $ppt.ActivePresentation

# Slide index -> (word that used to start the second half, word that was
# already the whole second run). Both get concatenated into the new
# second run's text.
$targets = @(
    @{ Slide = 4;  First = "Token";     Second = "Tester" },   # Example: Token Tester -> Example: TokenTester
    @{ Slide = 6;  First = "Point";     Second = "Setter" },   # Example: Point Setter -> Example: PointSetter
    @{ Slide = 10; First = "String";    Second = "Checker" },  # Example: String Checker -> Example: StringChecker
    @{ Slide = 13; First = "Reference"; Second = "Tester" },   # Example: Reference Tester -> Example: ReferenceTester
    @{ Slide = 20; First = "Equals";    Second = "Tester" },   # Example: Equals Tester -> Example: EqualsTester
    @{ Slide = 24; First = "Date";      Second = "Parser" },   # Exercise: Date Parser -> Exercise: DateParser
    @{ Slide = 25; First = "Box";       Second = "Volume" }    # Exercise: Box Volume -> Exercise: BoxVolume
)

foreach ($t in $targets) {
    $slide = $p.Slides.Item($t.Slide)
    $shape = $slide.Shapes.Item(1)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $camel = $t.First + $t.Second

    # The span to replace is "<First><space><Second>" at the end of the
    # title text; it becomes just "<First><Second>" (camelCase), joined
    # into the run that used to hold only "<Second>".
    $oldTail = $t.First + " " + $t.Second
    $start = $full.Length - $oldTail.Length + 1
    $len = $oldTail.Length

    $tr.Characters($start, $len).Text = $camel
}
